$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.308.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.642.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.69'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.651'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.47%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.83'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.389'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.90'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.117.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.133.47'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.635.35'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.70'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.59'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.30'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.89'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.15%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.41'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '586.54'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.60'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.05'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.60'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.414'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.12'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.14'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.07'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '159.71'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.03'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.46'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0604'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.104'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.636'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0255'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.58%  '
